$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B column text values (word tokens)
$ws.Range("B4").Value = "<that>"
$ws.Range("B9").Value = "<down>"
$ws.Range("B10").Value = "<alt>"
$ws.Range("B14").Value = "<they>"

# Update C column numeric values
$ws.Range("C2").Value = 6
$ws.Range("C3").Value = 16
$ws.Range("C4").Value = 11
$ws.Range("C6").Value = 14
$ws.Range("C7").Value = 12
$ws.Range("C8").Value = 16
$ws.Range("C9").Value = 14
$ws.Range("C10").Value = 8
$ws.Range("C11").Value = 12
$ws.Range("C12").Value = 6
$ws.Range("C13").Value = 12
$ws.Range("C14").Value = 10
$ws.Range("C15").Value = 15
$ws.Range("C16").Value = 13
$ws.Range("C17").Value = 17
$ws.Range("C18").Value = 4
